$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'30.469.22"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = "'1.867.54"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'235.52"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.07%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = "'0.4824"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = "'0.2804"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('D9').Value = "'0.06511"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('D10').Value = "'1.842.52"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.72%  '
$ws.Range('D11').Value = "'0.07436"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.16%  '
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').Value = "'5.069"
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Value = "'87.34"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.75%  '
$ws.Range('D15').Value = "'0.6469"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.60%  '
$ws.Range('D16').Value = "'30.448.23"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('E18').Value = '  -2.50%  '
$ws.Range('D19').Value = "'233.54"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.47%  '
$ws.Range('D20').Value = "'0.000007537"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.28%  '
$ws.Range('D21').Value = "'2.109.50"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').Value = "'1.000"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').Value = "'5.155"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.66%  '
$ws.Range('D24').Value = "'6.102"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.47%  '
$ws.Range('D25').Value = "'9.344"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.02%  '
$ws.Range('D26').Value = "'167.45"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.49%  '
$ws.Range('D27').Value = "'18.39"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('D28').Value = "'1.925"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.76%  '
$ws.Range('D29').Value = "'0.1029"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.42%  '
$ws.Range('D30').Value = "'1.373"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.58%  '
$ws.Range('D31').Value = "'4.275"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.70%  '
$ws.Range('D32').Value = "'4.014"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.16%  '
$ws.Range('D33').Value = "'0.04983"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.07%  '
$ws.Range('D34').Value = "'1.181"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.14%  '
$ws.Range('D35').Value = "'0.7463"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('D36').Value = "'1.0000"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('E37').Value = '  +0.40%  '
$ws.Range('E38').Value = '  +5.14%  '
$ws.Range('D39').Value = "'2.637"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.69%  '
$ws.Range('D40').Value = "'0.9189"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.50%  '
$ws.Range('D41').Value = "'2.056"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.59%  '
$ws.Range('D42').Value = "'106.29"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.53%  '
$ws.Range('D43').Value = "'0.9959"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('D44').Value = "'0.4205"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.05%  '
$ws.Range('D45').Value = "'5.547"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.74%  '
$ws.Range('D46').Value = "'7.257"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.73%  '
$ws.Range('D47').Value = "'61.93"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.17%  '
$ws.Range('D48').Value = "'0.1233"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.24%  '
$ws.Range('D49').Value = "'8.867"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.69%  '
$ws.Range('D50').Value = "'1.442"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.90%  '
$ws.Range('D51').Value = "'33.65"
$ws.Range('D51').Style = 'Normal'
